# Insert a new row of data at row 39, shifting all existing rows (39-147) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 39 - this shifts row 39..147 down to 40..148
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new data point
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44414
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = 100112040
$ws.Cells.Item(39, 7).Value = "Cilantro"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 300
$ws.Cells.Item(39, 11).Value = 15000
$ws.Cells.Item(39, 12).Value = 16000
$ws.Cells.Item(39, 13).Value = 15500
$ws.Cells.Item(39, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(39, 15).Value = "Región Metropolitana"
$ws.Cells.Item(39, 16).Value = 431
$ws.Cells.Item(39, 17).Value = 36
$ws.Cells.Item(39, 18).Value = "Hortaliza"

# Ensure the Fecha (date) cell uses the same number format as the rest of column D
$ws.Cells.Item(39, 4).NumberFormat = $ws.Cells.Item(40, 4).NumberFormat
